$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 113; this shifts existing rows 113-155
# down to 114-156 (data + formatting), matching the dimension change
# A1:R155 -> A1:R156 seen in the diff.
$ws.Rows.Item(113).Insert()

# Populate the newly inserted row 113 with the new record.
$ws.Range("A113").Value = 5
$ws.Range("B113").Value = "Macroferia Regional de Talca"
$ws.Range("C113").Value = "Maule"
$ws.Range("D113").Value2 = 44559
$ws.Range("E113").Value = 7
$ws.Range("F113").Value = 100112021
$ws.Range("G113").Value = "Ají"
$ws.Range("H113").Value = "Americana (o)"
$ws.Range("I113").Value = "Primera"
$ws.Range("J113").Value = 150
$ws.Range("K113").Value = 20000
$ws.Range("L113").Value = 20000
$ws.Range("M113").Value = 20000
$ws.Range("N113").Value = "`$/caja 14 kilos"
$ws.Range("O113").Value = "Región del Maule"
$ws.Range("P113").Value = 1429
$ws.Range("Q113").Value = 14
$ws.Range("R113").Value = "Hortaliza"
